$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user row (row 3) - a 3rd "personnel" record in the DB export
$ws.Range("A3").Value = "фамилия 3"
$ws.Range("B3").Value = "имя 3"
$ws.Range("C3").Value = "отчество 3"
$ws.Range("D3").Value = 1960
$ws.Range("E3").Value = 64
$ws.Range("F3").Value = 21
$ws.Range("G3").Value = 33
$ws.Range("H3").Value = "another@email.com"
$ws.Range("I3").Value = "qwerty1"
$ws.Range("J3").Value = "Мужской"
$ws.Range("K3").Value = "В браке"
$ws.Range("L3").Value = 2

# Hyperlink the e-mail cell, like the two rows above it
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:another@email.com") | Out-Null

# Re-apply the same cell format the other hyperlinked e-mail cells use
# (copy H1's format onto H3 so it reuses the existing "Hyperlink" style)
$ws.Range("H1").Copy()
$ws.Range("H3").PasteSpecial(-4122)

# Move the active selection onto the newly entered cell
$ws.Range("L3").Select() | Out-Null
